$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1086
$ws1.Range("F4").Value = 1668
$ws1.Range("F5").Value = 753
$ws1.Range("F6").Value = 172

# Sheet "全部类型": update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1086
$ws4.Range("F4").Value = 1668
$ws4.Range("F6").Value = 753
$ws4.Range("F7").Value = 172
